$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the Fitness column (C) for rows 189 through 252 (Generation 187-250)
# from 7310 to 7293.
$ws.Range("C189:C252").Value = 7293
